# Renaming 3,2,1 ph sheets to LLL, LL, LG (per commit message), plus
# the associated view/selection state and one new formatted cell that
# were captured in the same save.

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets -------------------------------------------------
$renames = @{
    "3ph_max"       = "LLL_max"
    "3ph_min"       = "LLL_min"
    "2ph_max"       = "LL_max"
    "2ph_min"       = "LL_min"
    "1ph_max"       = "LG_max"
    "1ph_min"       = "LG_min"
    "3ph_max_fault" = "LLL_max_fault"
    "3ph_min_fault" = "LLL_min_fault"
    "2ph_max_fault" = "LL_max_fault"
    "2ph_min_fault" = "LL_min_fault"
    "1ph_max_fault" = "LG_max_fault"
    "1ph_min_fault" = "LG_min_fault"
}

foreach ($oldName in $renames.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renames[$oldName]
}

# --- 2. Add a new bold, otherwise-empty cell on LLG_max (row 27) ------
$wsLLGmax = $wb.Worksheets.Item("LLG_max")
$wsLLGmax.Range("D27").Font.Bold = $true

# --- 3. Restore/update the various sheet selections recorded in the
#        saved view state. The last sheet selected below ends up the
#        active (tabSelected) sheet, matching "LL_max_fault".
$wb.Worksheets.Item("LG_min").Range("F1:F6").Select()
$wb.Worksheets.Item("LG_min_fault").Range("L14").Select()
$wsLLGmax.Range("D30").Select()
$wb.Worksheets.Item("LLG_max_fault").Range("D29").Select()
$wb.Worksheets.Item("LLG_min_fault").Range("O26").Select()
$wb.Worksheets.Item("LL_max_fault").Range("J13").Select()
